$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9-17 down to 10-18.
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the new record (values copied from the row that
# used to be there, date/price/origin updated per the new weekly entry).
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(9, 3).Value = "Coquimbo"
$ws.Cells.Item(9, 4).Value = 44721
$ws.Cells.Item(9, 4).NumberFormat = $ws.Cells.Item(10, 4).NumberFormat
$ws.Cells.Item(9, 5).Value = 4
$ws.Cells.Item(9, 6).Value = 100112026
$ws.Cells.Item(9, 7).Value = "Haba"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 500
$ws.Cells.Item(9, 11).Value = 14500
$ws.Cells.Item(9, 12).Value = 15000
$ws.Cells.Item(9, 13).Value = 14750
$ws.Cells.Item(9, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(9, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(9, 16).Value = 590
$ws.Cells.Item(9, 17).Value = 25
$ws.Cells.Item(9, 18).Value = "Hortaliza"
